$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Kayitlar")

$rng = $ws.Range("A1:G2")

# Every value in this sheet is stored as plain text (numbers/dates
# included), matching the rest of the workbook's sheets. Forcing a text
# number format while writing prevents auto-coercion into numeric/date
# values, then clearing formats afterwards drops back to the workbook's
# default (unstyled) cell style.
$rng.NumberFormat = "@"

# Header row
$ws.Range("A1").Value = "Kayıt No"
$ws.Range("B1").Value = "Tarih"
$ws.Range("C1").Value = "Birim"
$ws.Range("D1").Value = "Dosya Sayısı"
$ws.Range("E1").Value = "Parsel Sayısı"
$ws.Range("F1").Value = "İş"
$ws.Range("G1").Value = "Personeller"

# Data row - new record (Kayıt No 2)
$ws.Range("A2").Value = "2"
$ws.Range("B2").Value = "2025-07-16"
$ws.Range("C2").Value = "İlçe"
$ws.Range("D2").Value = "2"
$ws.Range("E2").Value = "3"
$ws.Range("F2").Value = "İfraz"
$ws.Range("G2").Value = "Gökhan ELGÜL"

$rng.ClearFormats()
